$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header: G1 "FT_Goals_H" -> "FT_Goals_A"
$ws.Range("G1").Value = "FT_Goals_A"

# Append 12 new match rows (312-323) for E1 14/01/2023
# Row 312
$ws.Range("A312").Value = "E1"
$ws.Range("B312").Value = "14/01/2023"
$ws.Range("C312").Value = "12:30"
$ws.Range("D312").Value = "Rotherham"
$ws.Range("E312").Value = "Blackburn"
$ws.Range("F312").Value = 4
$ws.Range("G312").Value = 0
$ws.Range("H312").Value = "H"
$ws.Range("I312").Value = 1
$ws.Range("J312").Value = 0
$ws.Range("K312").Value = "H"
$ws.Range("L312").Value = "D Whitestone"
$ws.Range("M312").Value = 8
$ws.Range("N312").Value = 9
$ws.Range("O312").Value = 4
$ws.Range("P312").Value = 3
$ws.Range("Q312").Value = 7
$ws.Range("R312").Value = 1
$ws.Range("S312").Value = 4
$ws.Range("T312").Value = 2
$ws.Range("U312").Value = 1
$ws.Range("V312").Value = 0
$ws.Range("W312").Value = 0
$ws.Range("X312").Value = 0
$ws.Range("Y312").Value = 3.4
$ws.Range("Z312").Value = 3.4
$ws.Range("AA312").Value = 2.15
$ws.Range("AB312").Value = 2.1
$ws.Range("AC312").Value = 1.73

# Row 313
$ws.Range("A313").Value = "E1"
$ws.Range("B313").Value = "14/01/2023"
$ws.Range("C313").Value = "15:00"
$ws.Range("D313").Value = "Bristol City"
$ws.Range("E313").Value = "Birmingham"
$ws.Range("F313").Value = 4
$ws.Range("G313").Value = 2
$ws.Range("H313").Value = "H"
$ws.Range("I313").Value = 2
$ws.Range("J313").Value = 1
$ws.Range("K313").Value = "H"
$ws.Range("L313").Value = "J Busby"
$ws.Range("M313").Value = 10
$ws.Range("N313").Value = 14
$ws.Range("O313").Value = 5
$ws.Range("P313").Value = 4
$ws.Range("Q313").Value = 14
$ws.Range("R313").Value = 6
$ws.Range("S313").Value = 5
$ws.Range("T313").Value = 2
$ws.Range("U313").Value = 2
$ws.Range("V313").Value = 1
$ws.Range("W313").Value = 0
$ws.Range("X313").Value = 0
$ws.Range("Y313").Value = 2.25
$ws.Range("Z313").Value = 3.3
$ws.Range("AA313").Value = 3.3
$ws.Range("AB313").Value = 2.1
$ws.Range("AC313").Value = 1.73

# Row 314
$ws.Range("A314").Value = "E1"
$ws.Range("B314").Value = "14/01/2023"
$ws.Range("C314").Value = "15:00"
$ws.Range("D314").Value = "Burnley"
$ws.Range("E314").Value = "Coventry"
$ws.Range("F314").Value = 1
$ws.Range("G314").Value = 0
$ws.Range("H314").Value = "H"
$ws.Range("I314").Value = 0
$ws.Range("J314").Value = 0
$ws.Range("K314").Value = "D"
$ws.Range("L314").Value = "D Webb"
$ws.Range("M314").Value = 10
$ws.Range("N314").Value = 8
$ws.Range("O314").Value = 5
$ws.Range("P314").Value = 0
$ws.Range("Q314").Value = 4
$ws.Range("R314").Value = 8
$ws.Range("S314").Value = 3
$ws.Range("T314").Value = 5
$ws.Range("U314").Value = 3
$ws.Range("V314").Value = 3
$ws.Range("W314").Value = 0
$ws.Range("X314").Value = 0
$ws.Range("Y314").Value = 1.7
$ws.Range("Z314").Value = 3.8
$ws.Range("AA314").Value = 5
$ws.Range("AB314").Value = 1.93
$ws.Range("AC314").Value = 1.93

# Row 315
$ws.Range("A315").Value = "E1"
$ws.Range("B315").Value = "14/01/2023"
$ws.Range("C315").Value = "15:00"
$ws.Range("D315").Value = "Cardiff"
$ws.Range("E315").Value = "Wigan"
$ws.Range("F315").Value = 1
$ws.Range("G315").Value = 1
$ws.Range("H315").Value = "D"
$ws.Range("I315").Value = 0
$ws.Range("J315").Value = 0
$ws.Range("K315").Value = "D"
$ws.Range("L315").Value = "J Simpson"
$ws.Range("M315").Value = 11
$ws.Range("N315").Value = 16
$ws.Range("O315").Value = 2
$ws.Range("P315").Value = 4
$ws.Range("Q315").Value = 7
$ws.Range("R315").Value = 11
$ws.Range("S315").Value = 3
$ws.Range("T315").Value = 3
$ws.Range("U315").Value = 4
$ws.Range("V315").Value = 2
$ws.Range("W315").Value = 0
$ws.Range("X315").Value = 0
$ws.Range("Y315").Value = 1.95
$ws.Range("Z315").Value = 3.4
$ws.Range("AA315").Value = 4
$ws.Range("AB315").Value = 2.3
$ws.Range("AC315").Value = 1.62

# Row 316
$ws.Range("A316").Value = "E1"
$ws.Range("B316").Value = "14/01/2023"
$ws.Range("C316").Value = "15:00"
$ws.Range("D316").Value = "Hull"
$ws.Range("E316").Value = "Huddersfield"
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 1
$ws.Range("H316").Value = "D"
$ws.Range("I316").Value = 0
$ws.Range("J316").Value = 1
$ws.Range("K316").Value = "A"
$ws.Range("L316").Value = "T Bramall"
$ws.Range("M316").Value = 14
$ws.Range("N316").Value = 8
$ws.Range("O316").Value = 5
$ws.Range("P316").Value = 6
$ws.Range("Q316").Value = 4
$ws.Range("R316").Value = 5
$ws.Range("S316").Value = 4
$ws.Range("T316").Value = 4
$ws.Range("U316").Value = 0
$ws.Range("V316").Value = 2
$ws.Range("W316").Value = 0
$ws.Range("X316").Value = 0
$ws.Range("Y316").Value = 1.95
$ws.Range("Z316").Value = 3.4
$ws.Range("AA316").Value = 4.2
$ws.Range("AB316").Value = 2.3
$ws.Range("AC316").Value = 1.62

# Row 317
$ws.Range("A317").Value = "E1"
$ws.Range("B317").Value = "14/01/2023"
$ws.Range("C317").Value = "15:00"
$ws.Range("D317").Value = "Luton"
$ws.Range("E317").Value = "West Brom"
$ws.Range("F317").Value = 2
$ws.Range("G317").Value = 3
$ws.Range("H317").Value = "A"
$ws.Range("I317").Value = 2
$ws.Range("J317").Value = 1
$ws.Range("K317").Value = "H"
$ws.Range("L317").Value = "J Smith"
$ws.Range("M317").Value = 10
$ws.Range("N317").Value = 18
$ws.Range("O317").Value = 4
$ws.Range("P317").Value = 7
$ws.Range("Q317").Value = 15
$ws.Range("R317").Value = 9
$ws.Range("S317").Value = 2
$ws.Range("T317").Value = 9
$ws.Range("U317").Value = 3
$ws.Range("V317").Value = 2
$ws.Range("W317").Value = 0
$ws.Range("X317").Value = 0
$ws.Range("Y317").Value = 3
$ws.Range("Z317").Value = 3.2
$ws.Range("AA317").Value = 2.45
$ws.Range("AB317").Value = 2.2
$ws.Range("AC317").Value = 1.67

# Row 318
$ws.Range("A318").Value = "E1"
$ws.Range("B318").Value = "14/01/2023"
$ws.Range("C318").Value = "15:00"
$ws.Range("D318").Value = "Middlesbrough"
$ws.Range("E318").Value = "Millwall"
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 0
$ws.Range("H318").Value = "H"
$ws.Range("I318").Value = 0
$ws.Range("J318").Value = 0
$ws.Range("K318").Value = "D"
$ws.Range("L318").Value = "M Donohue"
$ws.Range("M318").Value = 7
$ws.Range("N318").Value = 8
$ws.Range("O318").Value = 3
$ws.Range("P318").Value = 2
$ws.Range("Q318").Value = 8
$ws.Range("R318").Value = 13
$ws.Range("S318").Value = 6
$ws.Range("T318").Value = 3
$ws.Range("U318").Value = 2
$ws.Range("V318").Value = 3
$ws.Range("W318").Value = 0
$ws.Range("X318").Value = 0
$ws.Range("Y318").Value = 1.95
$ws.Range("Z318").Value = 3.5
$ws.Range("AA318").Value = 4
$ws.Range("AB318").Value = 2.1
$ws.Range("AC318").Value = 1.73

# Row 319
$ws.Range("A319").Value = "E1"
$ws.Range("B319").Value = "14/01/2023"
$ws.Range("C319").Value = "15:00"
$ws.Range("D319").Value = "Preston"
$ws.Range("E319").Value = "Norwich"
$ws.Range("F319").Value = 0
$ws.Range("G319").Value = 4
$ws.Range("H319").Value = "A"
$ws.Range("I319").Value = 0
$ws.Range("J319").Value = 3
$ws.Range("K319").Value = "A"
$ws.Range("L319").Value = "S Martin"
$ws.Range("M319").Value = 9
$ws.Range("N319").Value = 22
$ws.Range("O319").Value = 2
$ws.Range("P319").Value = 8
$ws.Range("Q319").Value = 17
$ws.Range("R319").Value = 11
$ws.Range("S319").Value = 4
$ws.Range("T319").Value = 7
$ws.Range("U319").Value = 4
$ws.Range("V319").Value = 4
$ws.Range("W319").Value = 0
$ws.Range("X319").Value = 0
$ws.Range("Y319").Value = 2.88
$ws.Range("Z319").Value = 3.3
$ws.Range("AA319").Value = 2.5
$ws.Range("AB319").Value = 2.1
$ws.Range("AC319").Value = 1.73

# Row 320
$ws.Range("A320").Value = "E1"
$ws.Range("B320").Value = "14/01/2023"
$ws.Range("C320").Value = "15:00"
$ws.Range("D320").Value = "Reading"
$ws.Range("E320").Value = "QPR"
$ws.Range("F320").Value = 2
$ws.Range("G320").Value = 2
$ws.Range("H320").Value = "D"
$ws.Range("I320").Value = 2
$ws.Range("J320").Value = 0
$ws.Range("K320").Value = "H"
$ws.Range("L320").Value = "D Bond"
$ws.Range("M320").Value = 10
$ws.Range("N320").Value = 17
$ws.Range("O320").Value = 3
$ws.Range("P320").Value = 4
$ws.Range("Q320").Value = 11
$ws.Range("R320").Value = 9
$ws.Range("S320").Value = 5
$ws.Range("T320").Value = 7
$ws.Range("U320").Value = 2
$ws.Range("V320").Value = 2
$ws.Range("W320").Value = 0
$ws.Range("X320").Value = 0
$ws.Range("Y320").Value = 2.7
$ws.Range("Z320").Value = 3.2
$ws.Range("AA320").Value = 2.75
$ws.Range("AB320").Value = 2.2
$ws.Range("AC320").Value = 1.67

# Row 321
$ws.Range("A321").Value = "E1"
$ws.Range("B321").Value = "14/01/2023"
$ws.Range("C321").Value = "15:00"
$ws.Range("D321").Value = "Sheffield United"
$ws.Range("E321").Value = "Stoke"
$ws.Range("F321").Value = 3
$ws.Range("G321").Value = 1
$ws.Range("H321").Value = "H"
$ws.Range("I321").Value = 2
$ws.Range("J321").Value = 1
$ws.Range("K321").Value = "H"
$ws.Range("L321").Value = "M Salisbury"
$ws.Range("M321").Value = 8
$ws.Range("N321").Value = 9
$ws.Range("O321").Value = 3
$ws.Range("P321").Value = 5
$ws.Range("Q321").Value = 7
$ws.Range("R321").Value = 13
$ws.Range("S321").Value = 2
$ws.Range("T321").Value = 6
$ws.Range("U321").Value = 2
$ws.Range("V321").Value = 2
$ws.Range("W321").Value = 0
$ws.Range("X321").Value = 0
$ws.Range("Y321").Value = 1.73
$ws.Range("Z321").Value = 3.75
$ws.Range("AA321").Value = 5
$ws.Range("AB321").Value = 2.02
$ws.Range("AC321").Value = 1.83

# Row 322
$ws.Range("A322").Value = "E1"
$ws.Range("B322").Value = "14/01/2023"
$ws.Range("C322").Value = "15:00"
$ws.Range("D322").Value = "Sunderland"
$ws.Range("E322").Value = "Swansea"
$ws.Range("F322").Value = 1
$ws.Range("G322").Value = 3
$ws.Range("H322").Value = "A"
$ws.Range("I322").Value = 0
$ws.Range("J322").Value = 0
$ws.Range("K322").Value = "D"
$ws.Range("L322").Value = "K Stroud"
$ws.Range("M322").Value = 8
$ws.Range("N322").Value = 18
$ws.Range("O322").Value = 3
$ws.Range("P322").Value = 6
$ws.Range("Q322").Value = 15
$ws.Range("R322").Value = 15
$ws.Range("S322").Value = 1
$ws.Range("T322").Value = 4
$ws.Range("U322").Value = 1
$ws.Range("V322").Value = 4
$ws.Range("W322").Value = 1
$ws.Range("X322").Value = 0
$ws.Range("Y322").Value = 2.6
$ws.Range("Z322").Value = 3.3
$ws.Range("AA322").Value = 2.75
$ws.Range("AB322").Value = 2
$ws.Range("AC322").Value = 1.85

# Row 323
$ws.Range("A323").Value = "E1"
$ws.Range("B323").Value = "14/01/2023"
$ws.Range("C323").Value = "15:00"
$ws.Range("D323").Value = "Watford"
$ws.Range("E323").Value = "Blackpool"
$ws.Range("F323").Value = 2
$ws.Range("G323").Value = 0
$ws.Range("H323").Value = "H"
$ws.Range("I323").Value = 0
$ws.Range("J323").Value = 0
$ws.Range("K323").Value = "D"
$ws.Range("L323").Value = "G Ward"
$ws.Range("M323").Value = 11
$ws.Range("N323").Value = 6
$ws.Range("O323").Value = 5
$ws.Range("P323").Value = 2
$ws.Range("Q323").Value = 11
$ws.Range("R323").Value = 14
$ws.Range("S323").Value = 2
$ws.Range("T323").Value = 3
$ws.Range("U323").Value = 2
$ws.Range("V323").Value = 3
$ws.Range("W323").Value = 0
$ws.Range("X323").Value = 0
$ws.Range("Y323").Value = 1.8
$ws.Range("Z323").Value = 3.75
$ws.Range("AA323").Value = 4.33
$ws.Range("AB323").Value = 1.95
$ws.Range("AC323").Value = 1.9

